$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert two new bulleted list items ("Wedges can draw" and
#    "Game ends when leaders are destroyed") right before the
#    "After Fully Playable" paragraph. We clone the formatting of an
#    existing top-level (ilvl=0) ListParagraph bullet so the new
#    paragraphs pick up the correct style / numbering / fonts, then
#    overwrite their text.
# ------------------------------------------------------------------

function Find-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text.TrimEnd("`r", "`n", "`x07") -eq $text) {
            return $para
        }
    }
    return $null
}

$templatePara = Find-ParagraphByText $d "Abilities need to be fully working"
$templateRange = $templatePara.Range.Duplicate

$newItems = @("Wedges can draw", "Game ends when leaders are destroyed")

foreach ($itemText in $newItems) {
    $targetPara = Find-ParagraphByText $d "After Fully Playable"
    $insertAt = $targetPara.Range.Start
    $insertRange = $d.Range($insertAt, $insertAt)
    $insertRange.FormattedText = $templateRange.FormattedText

    # re-resolve "After Fully Playable" since the document shifted; the
    # paragraph immediately before it is the one we just inserted.
    $freshTargetPara = Find-ParagraphByText $d "After Fully Playable"
    $newPara = $d.Paragraphs.Item($freshTargetPara.Index - 1)
    $newTextRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
    $newTextRange.Text = $itemText
}

# ------------------------------------------------------------------
# 2. Remove the "_GoBack" bookmark from the "After Fully Playable"
#    paragraph (it gets relocated in step 3).
# ------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 3. Split "Sounds to go with animations (attack, heal, destroy,
#    etc.)" after "Sounds to " and re-insert the "_GoBack" bookmark
#    at that split point.
# ------------------------------------------------------------------

$soundsPara = Find-ParagraphByText $d "Sounds to go with animations (attack, heal, destroy, etc.)"
$splitPos = $soundsPara.Range.Start + "Sounds to ".Length
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
